$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.273.23"
$ws.Range("E2").Value = "  -3.41%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.161.27"
$ws.Range("E3").Value = "  -2.03%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "607.89"
$ws.Range("E5").Value = "  +0.84%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.09"
$ws.Range("E6").Value = "  -6.05%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.159.22"
$ws.Range("E8").Value = "  -2.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.526"
$ws.Range("E9").Value = "  -3.30%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.151"
$ws.Range("E10").Value = "  -6.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.51"
$ws.Range("E11").Value = "  -4.10%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.477"
$ws.Range("E12").Value = "  -4.81%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000257"
$ws.Range("E13").Value = "  -4.61%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.25"
$ws.Range("E14").Value = "  -6.80%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.673.44"
$ws.Range("E15").Value = "  -2.31%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.262.30"
$ws.Range("E16").Value = "  -3.48%  "
$ws.Range("E17").Value = "  +1.02%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.151.77"
$ws.Range("E18").Value = "  -2.36%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.96"
$ws.Range("E19").Value = "  -4.24%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "479.37"
$ws.Range("E20").Value = "  -5.51%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.61"
$ws.Range("E21").Value = "  -4.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.711"
$ws.Range("E22").Value = "  -3.87%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.75"
$ws.Range("E23").Value = "  -3.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.79"
$ws.Range("E24").Value = "  -5.59%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.89"
$ws.Range("E25").Value = "  -3.52%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.89"
$ws.Range("E27").Value = "  -3.54%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.47"
$ws.Range("E28").Value = "  -6.45%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.21"
$ws.Range("E29").Value = "  -6.32%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.119"
$ws.Range("E30").Value = "  -27.96%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.83"
$ws.Range("E31").Value = "  -1.85%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.77"
$ws.Range("E32").Value = "  -4.98%  "
$ws.Range("E33").Value = "  -0.10%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.27"
$ws.Range("E34").Value = "  -6.97%  "
$ws.Range("E35").Value = "  -5.42%  "
$ws.Range("B36").Value = "Filecoin"
$ws.Range("C36").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.03"
$ws.Range("E36").Value = "  -5.51%  "
$ws.Range("B37").Value = "OKB"
$ws.Range("C37").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "54.18"
$ws.Range("E37").Value = "  -2.26%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0730"
$ws.Range("E38").Value = "  -7.69%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "451.96"
$ws.Range("E39").Value = "  -8.98%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.98"
$ws.Range("E40").Value = "  -6.62%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0398"
$ws.Range("E41").Value = "  -5.37%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.120"
$ws.Range("E42").Value = "  -6.50%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.47"
$ws.Range("E43").Value = "  -2.76%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.857.73"
$ws.Range("E44").Value = "  -2.90%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.270"
$ws.Range("E45").Value = "  -8.16%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.27"
$ws.Range("E46").Value = "  -7.64%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "26.42"
$ws.Range("E47").Value = "  -6.37%  "
$ws.Range("B48").Value = "USDe"
$ws.Range("C48").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.999"
$ws.Range("E48").Value = "  -0.03%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.31"
$ws.Range("E49").Value = "  -3.95%  "
$ws.Range("E50").Value = "  -3.10%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "117.82"
$ws.Range("E51").Value = "  -3.23%  "
